# Generate Report for Handback
# Update the generated timestamps on the handback status report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# "0500ae0c-...md" row - Latest HO Xliff Generate Date (Overview!G3)
# and the same value mirrored on de-de!H3 (Correspond Handoff Datetime)
$wsOverview.Range("G3").Value = "2016-08-22 06:44:31"
$wsDeDe.Range("H3").Value = "2016-08-22 06:44:31"

# zh-cn sheet, "0500ae0c-...md" row:
# Correspond Handoff Datetime (H3) and Correspond Handback DateTime (K3)
$wsZhCn.Range("H3").Value = "2016-08-22 06:44:27"
$wsZhCn.Range("K3").Value = "2016-08-22 06:44:45"

# de-de sheet, "0500ae0c-...md" row: Correspond Handback DateTime (K3)
$wsDeDe.Range("K3").Value = "2016-08-22 06:44:52"
